$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.071.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.357.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.51%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.74%  '
$ws.Range('E7').Value = '  +1.59%  '
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.628'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0934'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').Value = '  +10.71%  '
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.35'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.712.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.452.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +9.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.960.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('E22').Value = '  -1.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '252.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +10.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.88'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.89%  '
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.86'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.37%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.45%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0917'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.53%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.02%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.132'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.62%  '
$ws.Range('E37').Value = '  +2.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.59%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('E40').Value = '  +10.39%  '
$ws.Range('E41').Value = '  +14.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.66%  '
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.74%  '
$ws.Range('E47').Value = '  +9.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '110.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.04%  '
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0995'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.494.19'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.36%  '
